# Updates the "cryptos" price-ticker sheet: refreshed Price (D) and
# Volume(1h) (E) columns for the latest snapshot.
#
# Plain Range.Value assignment is used wherever the target text is safe
# (percentages always carry padding/%, and several prices contain more
# than one "." so Excel cannot mistake them for numbers). For prices that
# DO look like a plain number (e.g. "1.00", "526.29") a straight
# Range.Value assignment would get auto-converted to a numeric value by
# Excels normal typed-entry parsing, which is not what the source data
# contains (these are literal text strings). To avoid that - and to avoid
# leaving a NumberFormat/style override behind on the target cell (e.g. via
# a quote-prefix or a "@" text format) - such values are staged in a throw-
# away helper cell that is formatted as Text, then copied and pasted back
# with Paste Special > Values, which carries over the literal text without
# carrying over the helpers formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$ws.Range("D2").Value = '58.849.90'
$ws.Range("E2").Value = '  +2.53%  '

$ws.Range("D3").Value = '2.504.42'
$ws.Range("E3").Value = '  +3.50%  '

$ws.Range("E4").Value = '  +0.47%  '

$helper.Value = '526.29'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +2.71%  '

$helper.Value = '135.10'
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +4.94%  '

$helper.Value = '1.00'
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  +0.11%  '

$helper.Value = '0.567'
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +3.36%  '

$ws.Range("D9").Value = '2.537.94'
$ws.Range("E9").Value = '  +4.58%  '

$ws.Range("E10").Value = '  +3.78%  '

$ws.Range("E11").Value = '  -0.80%  '

$ws.Range("E12").Value = '  +0.52%  '

$ws.Range("E13").Value = '  +1.37%  '

$ws.Range("D14").Value = '2.996.26'
$ws.Range("E14").Value = '  +5.24%  '

$ws.Range("D15").Value = '58.799.10'
$ws.Range("E15").Value = '  +2.55%  '

$helper.Value = '22.41'
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  +4.64%  '

$ws.Range("E17").Value = '  +3.52%  '

$ws.Range("D18").Value = '2.542.77'
$ws.Range("E18").Value = '  +5.09%  '

$ws.Range("E19").Value = '  +3.85%  '

$helper.Value = '324.07'
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +3.19%  '

$ws.Range("E21").Value = '  +3.19%  '

$ws.Range("E22").Value = '  +8.29%  '

$helper.Value = '0.999'
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -0.03%  '

$helper.Value = '65.02'
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +2.37%  '

$helper.Value = '0.411'
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +2.47%  '

$ws.Range("E27").Value = '  +1.53%  '

$helper.Value = '7.54'
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  +4.82%  '

$ws.Range("D29").Value = '0.0₃0760'
$ws.Range("E29").Value = '  +6.15%  '

$ws.Range("E30").Value = '  +7.26%  '

$ws.Range("E31").Value = '  +4.94%  '

$helper.Value = '169.62'
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +0.22%  '

$helper.Value = '6.36'
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  +3.22%  '

$helper.Value = '0.998'
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)

$helper.Value = '0.992'
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -0.54%  '

$helper.Value = '18.29'
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +3.46%  '

$ws.Range("E37").Value = '  +0.49%  '

$ws.Range("E38").Value = '  +3.91%  '

$ws.Range("E39").Value = '  +5.73%  '

$helper.Value = '36.73'
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +1.15%  '

$helper.Value = '0.787'
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  +2.81%  '

$helper.Value = '280.56'
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +5.98%  '

$helper.Value = '134.80'
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +10.33%  '

$ws.Range("E44").Value = '  +3.91%  '

$helper.Value = '5.11'
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +5.01%  '

$helper.Value = '0.603'
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +3.68%  '

$ws.Range("E47").Value = '  +2.60%  '

$helper.Value = '0.0506'
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +5.50%  '

$ws.Range("E49").Value = '  +3.89%  '

$helper.Value = '17.19'
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +4.36%  '

$ws.Range("D51").Value = '1.759.95'
$ws.Range("E51").Value = '  +4.13%  '

$helper.Clear()
$excel.CutCopyMode = $false
